# Update 13C-MFA files (run and result) for SC and IO under WT-batch and chemostats
#
# - FluxData: insert a new flux row "EX_glc__D_e.f" right after "BIOMASS.f"
#   (pushing every following flux row down by one), update BIOMASS.f's
#   error value, and give the (now-shifted) EX_c5sugal_e.f row its own
#   value/error.
# - View state: FluxData becomes the active/selected sheet (instead of
#   MSData), all three sheets get their zoom bumped from 55% to 95%, and
#   each sheet's selection collapses down to a single top-left cell.

$wb = $excel.ActiveWorkbook

$msData = $wb.Worksheets.Item("MSData")
$fluxData = $wb.Worksheets.Item("FluxData")
$tracerData = $wb.Worksheets.Item("TracerData")

# --- FluxData: insert the new EX_glc__D_e.f flux row ------------------
# Before:
#   row2 BIOMASS.f        0.34           1E-06
#   row3 EX_c5sugal_e.f   0.0465512...   0.0046551...
# After:
#   row2 BIOMASS.f        0.34           0.0001
#   row3 EX_glc__D_e.f    6.43851654631153  3.14620763991209   <- new row
#   row4 EX_c5sugal_e.f   0.017361111111111 0.046551267319139
$fluxData.Rows.Item(3).Insert()

$fluxData.Range("A3").Value = "EX_glc__D_e.f"
$fluxData.Range("B3").Value = 6.43851654631153
$fluxData.Range("C3").Value = 3.14620763991209

$fluxData.Range("C2").Value = 0.0001

$fluxData.Range("B4").Value = 0.017361111111111
$fluxData.Range("C4").Value = 0.046551267319139

# --- View state ---------------------------------------------------------
# Every sheet's zoom goes from 55 to 95, and each selection collapses to
# a single cell. FluxData becomes the active sheet/tab (was MSData).

$msData.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 95
$msData.Range("A2").Select() | Out-Null

$fluxData.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 95
$fluxData.Range("A1").Select() | Out-Null

$tracerData.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 95
$tracerData.Range("A1").Select() | Out-Null

$fluxData.Activate() | Out-Null
